$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-45 down to 40-46
$ws.Rows(39).Insert()

# Copy formatting (date style) from the row above (row 38's D cell) into new D39
$ws.Range("D38").Copy()
$ws.Range("D39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Populate the new row 39 with data
$ws.Cells.Item(39, 1).Value = 3
$ws.Cells.Item(39, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(39, 3).Value = "Coquimbo"
$ws.Cells.Item(39, 4).Value = 44736
$ws.Cells.Item(39, 5).Value = 5
$ws.Cells.Item(39, 6).Value = 100112035
$ws.Cells.Item(39, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 82
$ws.Cells.Item(39, 11).Value = 16000
$ws.Cells.Item(39, 12).Value = 17000
$ws.Cells.Item(39, 13).Value = 16488
$ws.Cells.Item(39, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(39, 16).Value = 1099
$ws.Cells.Item(39, 17).Value = 15
$ws.Cells.Item(39, 18).Value = "Hortaliza"
